$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 243, pushing the existing rows 243:278 down to 244:279.
$ws.Rows("243:243").Insert()

# Populate the newly inserted row 243 with the new data record.
$ws.Range("A243").Value = 3
$ws.Range("B243").Value = "Femacal de La Calera"
$ws.Range("C243").Value = "Coquimbo"
$ws.Range("D243").Value = 45131
$ws.Range("E243").Value = 5
$ws.Range("F243").Value = 100112026
$ws.Range("G243").Value = "Haba"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 80
$ws.Range("K243").Value = 14500
$ws.Range("L243").Value = 15000
$ws.Range("M243").Value = 14750
$ws.Range("N243").Value = "$/saco 25 kilos"
$ws.Range("O243").Value = "Provincia de Limarí"
$ws.Range("P243").Value = 590
$ws.Range("Q243").Value = 25
$ws.Range("R243").Value = "Hortaliza"
